$d = $word.ActiveDocument

# --- Change 1 --------------------------------------------------------
# Merge the three runs (separated by a "gramStart"/"gramEnd" proofErr
# pair) that made up the first body paragraph back into a single run,
# by replacing the whole sentence with itself (Word collapses the
# proofreading markers and run splits on a Find/Replace).
$old1 = "The natural frequency is observed when the system is left at any point and the oscillations due to the spring force is observed. (When P(t) is not actively oscillating the system) Thus, the load should be taken as zero, while the system is oscillating with this frequency. Then any  frequency that solves the governing equation are the natural frequencies of the system."
$new1 = "The natural frequency is observed when the system is left at any point and the oscillations due to the spring force is observed. (When P(t) is not actively oscillating the system) Thus, the load should be taken as zero, while the system is oscillating with this frequency. Then any  frequency that solves the governing equation are the natural frequencies of the system."
$d.Content.Find.Execute($old1, $false, $false, $false, $false, $false, $true, 1, $false, $new1, 2) | Out-Null

# --- Change 2 --------------------------------------------------------
# Fix the "bou|ndary" split (spellStart/spellEnd proofErr pair) so the
# word "boundary" and its surrounding sentence is a single run again.
$old2 = ". Using this boundary condition Equation 11 can be followed to find "
$new2 = ". Using this boundary condition Equation 11 can be followed to find "
$d.Content.Find.Execute($old2, $false, $false, $false, $false, $false, $true, 1, $false, $new2, 2) | Out-Null

# --- Change 3 --------------------------------------------------------
# Remove the leftover "_GoBack" bookmark (Word's "last edit position"
# marker) that trailed the oMathPara paragraph.
try {
    $goBack = $d.Bookmarks("_GoBack")
    if ($goBack -ne $null) {
        $goBack.Delete()
    }
} catch {
    # no "_GoBack" bookmark present - nothing to remove
}
